$wb = $excel.ActiveWorkbook

# --- Update the conversion text on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.51 = 25676.25 pesos`n✅ 25676.25 pesos = 6.49 = 966.35 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 153.7
$wsTasas.Range("O10").Value = 3946.44
$wsTasas.Range("N12").Value = 3958.99
$wsTasas.Range("O12").Value = 149
